$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets (groupBy, orderBy) after the last existing
#    sheet (groupCols).  Excel re-numbers sheetId's and flips the active tab
#    to the newly inserted (and now last/selected) sheet automatically.
# ---------------------------------------------------------------------------
$groupCols = $wb.Worksheets.Item($wb.Worksheets.Count)
$groupBy = $wb.Worksheets.Add($null, $groupCols)
$groupBy.Name = "groupBy"
$orderBy = $wb.Worksheets.Add($null, $groupBy)
$orderBy.Name = "orderBy"

$vertVert = $wb.Worksheets.Item("VertVert")
$groupColsSheet = $wb.Worksheets.Item("groupCols")

# ---------------------------------------------------------------------------
# 2. groupBy sheet - mirrors the simple "VertVert" layout but for a single
#    forEach with groupBy="divisionName".
# ---------------------------------------------------------------------------

# Row 1 - merged banner cell (copy style+value from VertVert!A1)
$vertVert.Range("A1:E1").Copy()
$groupBy.Range("A1:E1").PasteSpecial(-4122)
$groupBy.Range("A1").Value = '<jt:forEach items="${teams}" var="division" groupBy="divisionName">Division: ${division.obj.divisionName}'
$groupBy.Range("A1:E1").Merge()

# Row 2 - column headers (copy style from VertVert!A2:E2)
$vertVert.Range("A2:E2").Copy()
$groupBy.Range("A2:E2").PasteSpecial(-4122)
$groupBy.Range("A2").Value = "City"
$groupBy.Range("B2").Value = "Name"
$groupBy.Range("C2").Value = "Wins"
$groupBy.Range("D2").Value = "Losses"
$groupBy.Range("E2").Value = "Pct."

# Row 3 - data row (copy style from VertVert!A3:E3)
$vertVert.Range("A3:E3").Copy()
$groupBy.Range("A3:E3").PasteSpecial(-4122)
$groupBy.Range("A3").Value = '<jt:forEach items="${division.items}" var="team">${team.city}'
$groupBy.Range("B3").Value = '${team.name}'
$groupBy.Range("C3").Value = '${team.wins}'
$groupBy.Range("D3").Value = '${team.losses}'
$groupBy.Range("E3").Value = '${team.pct}'

$groupBy.Columns.Item(1).ColumnWidth = 14.2
$groupBy.Columns.Item(2).ColumnWidth = 14.6

$groupBy.Range("A1:E3").Select()

# ---------------------------------------------------------------------------
# 3. orderBy sheet - two independent blocks:
#    left block  (A:F) - forEach with orderBy="divisionName desc;pct"
#    right block (I:M) - forEach with groupBy + orderBy combined
# ---------------------------------------------------------------------------

# Left block header row (copy style from VertVert!A2:E2, which matches the
# plain bold/blue header look used for columns A:F here)
$vertVert.Range("A2:E2").Copy()
$orderBy.Range("A1:E1").PasteSpecial(-4122)
$vertVert.Range("A2").Copy()
$orderBy.Range("F1").PasteSpecial(-4122)

$orderBy.Range("A1").Value = "Division"
$orderBy.Range("B1").Value = "City"
$orderBy.Range("C1").Value = "Name"
$orderBy.Range("D1").Value = "Wins"
$orderBy.Range("E1").Value = "Losses"
$orderBy.Range("F1").Value = "Pct."

# Left block data row (copy style from VertVert!A3:E3 plus E3 for the pct cell)
$vertVert.Range("A3:E3").Copy()
$orderBy.Range("A2:E2").PasteSpecial(-4122)
$vertVert.Range("E3").Copy()
$orderBy.Range("F2").PasteSpecial(-4122)

$orderBy.Range("A2").Value = '<jt:forEach items="${teams}" var="team" orderBy="divisionName desc;pct">${team.divisionName}'
$orderBy.Range("B2").Value = '${team.city}'
$orderBy.Range("C2").Value = '${team.name}'
$orderBy.Range("D2").Value = '${team.wins}'
$orderBy.Range("E2").Value = '${team.losses}'
$orderBy.Range("F2").Value = '${team.pct}</jt:forEach>'

# Right block banner cell (copy style from groupCols!A1, the vertical merged
# banner style, but here laid out as a plain merged single row like groupBy)
$groupColsSheet.Range("A1:C1").Copy()
$orderBy.Range("I1:M1").PasteSpecial(-4122)
$orderBy.Range("I1").Value = '<jt:forEach items="${teams}" var="division" groupBy="divisionName" orderBy="divisionName desc;pct">Division: ${division.obj.divisionName}'
$orderBy.Range("I1:M1").Merge()

# Right block header row
$vertVert.Range("A2:E2").Copy()
$orderBy.Range("I2:M2").PasteSpecial(-4122)
$orderBy.Range("I2").Value = "City"
$orderBy.Range("J2").Value = "Name"
$orderBy.Range("K2").Value = "Wins"
$orderBy.Range("L2").Value = "Losses"
$orderBy.Range("M2").Value = "Pct."

# Right block data row
$vertVert.Range("A3:E3").Copy()
$orderBy.Range("I3:M3").PasteSpecial(-4122)
$orderBy.Range("I3").Value = '<jt:forEach items="${division.items}" var="team">${team.city}'
$orderBy.Range("J3").Value = '${team.name}'
$orderBy.Range("K3").Value = '${team.wins}'
$orderBy.Range("L3").Value = '${team.losses}'
$orderBy.Range("M3").Value = '${team.pct}'

$orderBy.Columns.Item(1).ColumnWidth = 14.2
$orderBy.Columns.Item(2).ColumnWidth = 14.2
$orderBy.Columns.Item(3).ColumnWidth = 14.6
$orderBy.Columns.Item(9).ColumnWidth = 14.2
$orderBy.Columns.Item(10).ColumnWidth = 14.6

$orderBy.Range("A1").Select()

# ---------------------------------------------------------------------------
# 4. "where" sheet tab should stay selected-less; just make sure orderBy is
#    the active tab like the real commit (last-added sheet stays active).
# ---------------------------------------------------------------------------
$orderBy.Activate()
